$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the fill formatting previously applied to A1 (resets it back to the
# default style, collapsing the cellXfs usage to the base style only).
$ws.Range("A1").Interior.Pattern = -4142

# Add the resume/summary text as a new cell C1 (new shared string entry).
$ws.Range("C1").Value = "Over three years of experience as a software quality analyst, worked on designing and maintaining scripts for automation testing. Knowledge of Maven, Jenkins, GIT, SQL & UNIX. Skilled in Java, Selenium WebDriver and Data-Driven Framework."

# Move the active selection to C3, matching the saved view state.
$ws.Range("C3").Select()
